$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handoff" - refresh the localization-status report with
# the newly generated source file (a316e2b1-...) replacing the previous run
# (3fec1a29-...), update the handoff timestamps, and clear out the
# now-stale handback info (no handback has happened yet for this new file).
# ---------------------------------------------------------------------------

$oldGuid = "3fec1a29-f8dc-499b-adbd-64b948868ea9"
$newGuid = "a316e2b1-e078-486b-b206-fde2de613608"

$oldZhXlf = "$oldGuid.f5863a18e1d03ab84f804234215253f1aa5fb653.zh-cn.xlf"
$newZhXlf = "$newGuid.7ab37ce20848ba4f81bb2258350f650138ce71f1.zh-cn.xlf"

$oldDeXlf = "$oldGuid.f5863a18e1d03ab84f804234215253f1aa5fb653.de-de.xlf"
$newDeXlf = "$newGuid.7ab37ce20848ba4f81bb2258350f650138ce71f1.de-de.xlf"

# ---------------------------------------------------------------------------
# "Overview" sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

$overviewLinkCell = $wsOverview.Range("B2")
$overviewLinkCell.Value = "e2e\$newGuid.md"
$overviewLinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/81ef94ae697e1e8d4be8e68ef2268af8e5437556/e2e/$newGuid.md"
$overviewLinkCell.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($overviewLinkCell, $overviewLinkUrl, "", "", "e2e\$newGuid.md")

$wsOverview.Range("G2").Value = "2016-08-21 17:05:42"

# ---------------------------------------------------------------------------
# "zh-cn" sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhCnLinkCell = $wsZhCn.Range("A2")
$zhCnLinkCell.Value = "$newGuid.md"
$zhCnLinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/81ef94ae697e1e8d4be8e68ef2268af8e5437556/e2e/$newGuid.md"
$zhCnLinkCell.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($zhCnLinkCell, $zhCnLinkUrl, "", "", "$newGuid.md")

$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = "2016-08-21 17:05:38"

# Target/handback files + handback datetime reset - no handback yet for the
# freshly generated handoff package.
$wsZhCn.Range("I2").Hyperlinks.Delete()
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# "de-de" sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deDeLinkCell = $wsDeDe.Range("A2")
$deDeLinkCell.Value = "$newGuid.md"
$deDeLinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/81ef94ae697e1e8d4be8e68ef2268af8e5437556/e2e/$newGuid.md"
$deDeLinkCell.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($deDeLinkCell, $deDeLinkUrl, "", "", "$newGuid.md")

$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = "2016-08-21 17:05:42"

$wsDeDe.Range("I2").Hyperlinks.Delete()
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Column width refresh for "Latest Target File" / "Latest Handback File" on
# the locale sheets now that those columns hold empty values.
# ---------------------------------------------------------------------------
$wsZhCn.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZhCn.Columns.Item(10).ColumnWidth = 21.7054770333426

$wsDeDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDeDe.Columns.Item(10).ColumnWidth = 21.7054770333426
